$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update destination text values
$ws.Range("B2").Value = "Dallas, Texas, United States of America"
$ws.Range("B3").Value = "Paris, France"

# Update Budget column from numeric to descriptive string buckets
$ws.Range("C2").Value = "$100 to $299"
$ws.Range("C3").Value = "More than $500"

# Update Rooms column from numeric to descriptive string values
$ws.Range("G2").Value = "2 Rooms"
$ws.Range("G3").Value = "1 Room"

# Update Pool column from boolean to Yes/No text
$ws.Range("I3").Value = "No"
$ws.Range("I2").Value = "Yes"

# Update date number format
$ws.Range("D2:E3").NumberFormat = "mm/dd/yyyy;@"

# Column widths for D:E
$ws.Range("D1:E1").ColumnWidth = 10.5546875

# Selection
$ws.Range("E5").Select()

# Page setup
$ws.PageSetup.Orientation = 1
